$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AO1").Value = 0.65147969523784
$ws.Range("BA1").Value = 0.81894404321217107
$ws.Range("BO1").Value = 0.96858229914515315
$ws.Range("AK2").Value = 0.97099443626253346
$ws.Range("BH2").Value = 0.79561549764251027
$ws.Range("AJ3").Value = 0.83526295817889462
$ws.Range("AN3").Value = 0.81761804893500389
$ws.Range("BN3").Value = 0.97672178218763395
$ws.Range("F5").Value = 0.76202548918266
$ws.Range("AA5").Value = 0.81495469111576291
$ws.Range("AD5").Value = 0.6108400052168772
$ws.Range("R6").Value = 0.91154672668052028
$ws.Range("U6").Value = 0.9123941204900845
$ws.Range("S7").Value = 0.99203712723317228
$ws.Range("W7").Value = 0.76062810809392201
$ws.Range("BF7").Value = 0.88520835366339612
$ws.Range("Y8").Value = 0.96074495695820361
$ws.Range("BF8").Value = 0.82003978833685109
$ws.Range("BP8").Value = 0.52586862018013747
$ws.Range("AG9").Value = 0.85072301590682187
$ws.Range("E11").Value = 0.71750111111920301
$ws.Range("M11").Value = 0.62633454855744142
$ws.Range("Y11").Value = 0.6577649778233825
$ws.Range("AY11").Value = 0.9401963837348819
$ws.Range("C12").Value = 0.9347811585505531
$ws.Range("J12").Value = 0.81609416882116292
$ws.Range("AZ12").Value = 0.54647400224859854
$ws.Range("A13").Value = 0.99359427332404615
$ws.Range("F13").Value = 0.95489281494336087
$ws.Range("BB13").Value = 0.87014401167411481
$ws.Range("BF13").Value = 0.71866594581182164
$ws.Range("BO13").Value = 0.90962422115509733
$ws.Range("E14").Value = 0.79347063789745576
$ws.Range("L14").Value = 0.65505414166005949
$ws.Range("AC14").Value = 0.85231960971033649
$ws.Range("W15").Value = 0.79596817045267509
$ws.Range("AD16").Value = 0.88036074832069988
$ws.Range("AF16").Value = 0.9187294158607926
$ws.Range("J17").Value = 0.96179617025379605
$ws.Range("X17").Value = 0.88850603788037374
$ws.Range("AM17").Value = 0.97507745371895804
$ws.Range("S18").Value = 0.87596711531190818
$ws.Range("AL18").Value = 0.82177290266556269
$ws.Range("BJ19").Value = 0.9185753430348037
$ws.Range("D20").Value = 0.98589826840531436
$ws.Range("U20").Value = 0.5465394155337131
$ws.Range("AI20").Value = 0.87159457553472741
$ws.Range("BE21").Value = 0.75007666583347388
$ws.Range("AY22").Value = 0.91008578190510359
$ws.Range("V23").Value = 0.65702515612594936
$ws.Range("AE23").Value = 0.75523966297103007
$ws.Range("BP23").Value = 0.65951544341062629
$ws.Range("BP24").Value = 0.82109131103772881
$ws.Range("M25").Value = 0.79668644565151592
$ws.Range("BK25").Value = 0.84092967408989172
$ws.Range("D26").Value = 0.77178468734555294
$ws.Range("J26").Value = 0.9129259046624798
$ws.Range("V26").Value = 0.78364078768371404
$ws.Range("BE26").Value = 0.71522893422978995
$ws.Range("U27").Value = 0.97444218680601891
$ws.Range("Y27").Value = 0.8225996716851609
$ws.Range("T28").Value = 0.85827055784407236
$ws.Range("AJ28").Value = 0.53534643265825754
$ws.Range("AF29").Value = 0.80621109056067963
$ws.Range("AP29").Value = 0.8110883784838766
$ws.Range("D30").Value = 0.68679424920294008
$ws.Range("O30").Value = 0.85413146806332085
$ws.Range("U30").Value = 0.75259700790765605
$ws.Range("AG31").Value = 0.65750173368098475
$ws.Range("AM31").Value = 0.76962253555935778
$ws.Range("AR32").Value = 0.70098087039547563
$ws.Range("Y33").Value = 0.88865440487097636
$ws.Range("AH33").Value = 0.89179809793034459
$ws.Range("BN33").Value = 0.70002505825814576
$ws.Range("A34").Value = 0.74779409741695813
$ws.Range("AS34").Value = 0.98343070180554237
$ws.Range("AV34").Value = 0.8790897873443646
$ws.Range("AV35").Value = 0.527186250455975
$ws.Range("BG35").Value = 0.97275160018262707
$ws.Range("BH35").Value = 0.91993322924845844
$ws.Range("Y36").Value = 0.88267137236149784
$ws.Range("AM36").Value = 0.79902850126479441
$ws.Range("BB36").Value = 0.60580900858828279
$ws.Range("AL37").Value = 0.83917706899696909
$ws.Range("N38").Value = 0.96771530895913171
$ws.Range("AN38").Value = 0.85477462924598879
$ws.Range("L39").Value = 0.91906171074099585
$ws.Range("AH39").Value = 0.76507733858966898
$ws.Range("Z40").Value = 0.63152801245107582
$ws.Range("BJ40").Value = 0.96316523842633628
$ws.Range("G41").Value = 0.98493602092014942
$ws.Range("AU41").Value = 0.81692300331438017
$ws.Range("AD42").Value = 0.92305963342774311
$ws.Range("S43").Value = 0.91801597979138516
$ws.Range("BJ43").Value = 0.91286049727244478
$ws.Range("F44").Value = 0.85725557874438429
$ws.Range("W44").Value = 0.98438988989991416
$ws.Range("BG44").Value = 0.75555073699327147
$ws.Range("AK45").Value = 0.94438301607796782
$ws.Range("BD45").Value = 0.7825118008108396
$ws.Range("B46").Value = 0.94350714091442323
$ws.Range("AV46").Value = 0.89933703110993424
$ws.Range("AI47").Value = 0.78959273623014303
$ws.Range("AP47").Value = 0.9877751705109481
$ws.Range("AS47").Value = 0.79653669469423516
$ws.Range("I48").Value = 0.90895316553191985
$ws.Range("BF48").Value = 0.93799987060607082
$ws.Range("AN49").Value = 0.83686390737544714
$ws.Range("AY49").Value = 0.75294180083773932
$ws.Range("A50").Value = 0.84209717340984924
$ws.Range("F51").Value = 0.68843690718313211
$ws.Range("AF51").Value = 0.82913285018096294
$ws.Range("AX51").Value = 0.81803122395986694
$ws.Range("Q52").Value = 0.99188088890230863
$ws.Range("AA52").Value = 0.85739621159766077
$ws.Range("L53").Value = 0.7671779866930708
$ws.Range("B55").Value = 0.75435203123047667
$ws.Range("F55").Value = 0.74426107451248902
$ws.Range("Z55").Value = 0.88587646769168027
$ws.Range("BG55").Value = 0.92727943714062366
$ws.Range("W56").Value = 0.87966448767620786
$ws.Range("BC57").Value = 0.58117016025875756
$ws.Range("AO60").Value = 0.8328790498878067
$ws.Range("BB60").Value = 0.71004296759584251
$ws.Range("BF60").Value = 0.83686269457448992
$ws.Range("AX61").Value = 0.67455373793628404
$ws.Range("AY62").Value = 0.85722568264934618
$ws.Range("R63").Value = 0.91783004457171691
$ws.Range("BI63").Value = 0.99562447703784351
$ws.Range("F64").Value = 0.68325617390370419
$ws.Range("BJ64").Value = 0.64393106539043266
$ws.Range("BM64").Value = 0.99513652253234697
$ws.Range("Q65").Value = 0.77963924324273504
$ws.Range("AB65").Value = 0.91877341003917345
$ws.Range("AT65").Value = 0.90928271624232637
$ws.Range("BM67").Value = 0.88136822795327918
